$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Q7)
$ws.Range("B9").Value = 0.4298197394642976
$ws.Range("C9").Value = 0.4298197394642976
$ws.Range("D9").Value = 0.2205190090366816
$ws.Range("E9").Value = 0.4695945155521747
$ws.Range("F9").Value = 0.2071926657105166
$ws.Range("G9").Value = 6

# Row 10 (Q8)
$ws.Range("B10").Value = 0.2695677751263204
$ws.Range("C10").Value = 0.2695677751263204
$ws.Range("D10").Value = 0.07820422228561637
$ws.Range("E10").Value = 0.2796501784115583
$ws.Range("F10").Value = 0.09113811139470071
$ws.Range("G10").Value = 3

# Row 11 (Q9)
$ws.Range("B11").Value = 0.1878533956106533
$ws.Range("C11").Value = 0.1878533956106533
$ws.Range("D11").Value = 0.03528889824245262
$ws.Range("E11").Value = 0.1878533956106533
$ws.Range("F11").ClearContents()
$ws.Range("G11").Value = 1
